# "Exclude the ServerData2 Folder" - append the new rows/values to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 gains a numeric value next to the existing "MSDP" label.
$ws.Range("B17").Value = 11

# New row 18 - a lone numeric value in column B.
$ws.Range("B18").Value = 111

# New block starting at row 21: a label/value pair followed by two more values.
$ws.Range("A21").Value = "HHH"
$ws.Range("B21").Value = "d"
$ws.Range("B22").Value = "ff"
$ws.Range("B23").Value = "ffff"

# Leave the selection where the user ended up after typing the last entry.
$ws.Range("B25").Select()
